# Add "last added films" retrieval support:
#  - insert a new "Date insertion" column before the existing "Zonedvd" column
#  - append a new "Date Sortie DVD" column at the end
#  - populate the new columns and update a few existing values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at H ("Zonedvd", "Rippé", "RIP Date", "Dvd Format" shift right by one)
$ws.Columns("H").Insert()

# New header for the inserted column, and a brand new trailing column M
$ws.Range("H1").Value = "Date insertion"
$ws.Range("M1").Value = "Date Sortie DVD"

# Force the cells below to be stored as plain text (dates / numbers as strings),
# matching how the rest of this text-only worksheet is laid out.
$dataCells = @(
    "H2","I2","J2","K2","L2","M2",
    "H3","I3","J3","K3","L3","M3",
    "H4","I4","J4","K4","L4","M4",
    "H5","I5","J5","K5","L5","M5",
    "H6","I6","J6","K6","L6","M6",
    "H7","I7","J7","K7","L7","M7",
    "H8","I8","J8","K8","L8","M8"
)
foreach ($addr in $dataCells) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C6").NumberFormat = "@"

# Row 2 - STANLEY KUBRICK / 2001 : L'ODYSSEE DE L'ESPACE
$ws.Range("E2").Value = "EN_SALLE"
$ws.Range("H2").Value = "10/01/2020"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = "31/07/0022"

# Row 3 - WONG KAR-WAI / 2046
$ws.Range("H3").Value = "10/01/2020"
$ws.Range("I3").Value = "1"
$ws.Range("J3").Value = "non"
$ws.Range("K3").Value = ""
$ws.Range("L3").Value = "DVD"
$ws.Range("M3").Value = ""

# Row 4 - ZACK SNYDER / 300
$ws.Range("H4").Value = "10/01/2020"
$ws.Range("I4").Value = "1"
$ws.Range("J4").Value = "oui"
$ws.Range("K4").Value = "06/05/2019"
$ws.Range("L4").Value = "DVD"
$ws.Range("M4").Value = ""

# Row 5 - JUDD APATOW / 40 ANS : MODE D'EMPLOI
$ws.Range("C5").Value = "2013"
$ws.Range("H5").Value = "10/01/2020"
$ws.Range("I5").Value = "1"
$ws.Range("J5").Value = "oui"
$ws.Range("K5").Value = "23/07/2019"
$ws.Range("L5").Value = "DVD"
$ws.Range("M5").Value = ""

# Row 6 - RAMIN BAHRANI / 99 HOMES
$ws.Range("C6").Value = "2016"
$ws.Range("H6").Value = "10/01/2020"
$ws.Range("I6").Value = "1"
$ws.Range("J6").Value = "non"
$ws.Range("K6").Value = ""
$ws.Range("L6").Value = "DVD"
$ws.Range("M6").Value = ""

# Row 7 - DAVID CRONENBERG / A HISTORY OF VIOLENCE
$ws.Range("H7").Value = "10/01/2020"
$ws.Range("I7").Value = "1"
$ws.Range("J7").Value = "non"
$ws.Range("K7").Value = ""
$ws.Range("L7").Value = "DVD"
$ws.Range("M7").Value = ""

# Row 8 - J. C. CHANDOR / A MOST VIOLENT YEAR
$ws.Range("H8").Value = "10/01/2020"
$ws.Range("I8").Value = "1"
$ws.Range("J8").Value = "non"
$ws.Range("K8").Value = ""
$ws.Range("L8").Value = "DVD"
$ws.Range("M8").Value = ""
